# Fixed naive component forecaster bug - Presentation state 11.02.
# The error-table rows shift down by one: each row's B:G values move to the
# next row down, and a brand new row 2 is populated with freshly computed
# (much smaller) error metrics, reflecting one additional matched quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 2 (the newly added/recomputed matched-error entry)
# (written in plain decimal form - the interpreter's numeric literal parser
# does not accept scientific "E" notation - but these round-trip to the
# exact same IEEE-754 double as the scientific-notation originals)
$ws.Range("B2").Value = 0.000000255470939168611611476704
$ws.Range("C2").Value = 0.000000801143367166097996470513
$ws.Range("D2").Value = 0.000000000003360900671017876935
$ws.Range("E2").Value = 0.000001833275939682261058719603
$ws.Range("F2").Value = 0.000001879105582334942090225725
$ws.Range("G2").Value = 15

# Row 3 (previously row 2's values)
$ws.Range("B3").Value = -0.005035358036557557
$ws.Range("C3").Value = 0.3637990656943072
$ws.Range("D3").Value = 0.1855351867275108
$ws.Range("E3").Value = 0.4307379559865961
$ws.Range("F3").Value = 0.4469673578304695
$ws.Range("G3").Value = 14

# Row 4 (previously row 3's values)
$ws.Range("B4").Value = -0.01898232632975465
$ws.Range("C4").Value = 0.3165751929851393
$ws.Range("D4").Value = 0.1363897707384557
$ws.Range("E4").Value = 0.3693098573534908
$ws.Range("F4").Value = 0.3838817902180699
$ws.Range("G4").Value = 13

# Row 5 (previously row 4's values)
$ws.Range("B5").Value = 0.00780077805212256
$ws.Range("C5").Value = 0.3794346132818944
$ws.Range("D5").Value = 0.1847968220473418
$ws.Range("E5").Value = 0.4298800088947401
$ws.Range("F5").Value = 0.4489210943938488
$ws.Range("G5").Value = 12

# Row 6 (previously row 5's values)
$ws.Range("B6").Value = 0.03532231998103826
$ws.Range("C6").Value = 0.2564852750112934
$ws.Range("D6").Value = 0.1004647918617584
$ws.Range("E6").Value = 0.3169618145167623
$ws.Range("F6").Value = 0.3303616777566779
$ws.Range("G6").Value = 11

# Row 7 (previously row 6's values)
$ws.Range("B7").Value = 0.002000156080238219
$ws.Range("C7").Value = 0.3449448151542904
$ws.Range("D7").Value = 0.1826303857083173
$ws.Range("E7").Value = 0.4273527649475516
$ws.Range("F7").Value = 0.4504644332784307
$ws.Range("G7").Value = 10

# Row 8 (previously row 7's values)
$ws.Range("B8").Value = -0.05875859174690282
$ws.Range("C8").Value = 0.3248866822092542
$ws.Range("D8").Value = 0.1283271032654252
$ws.Range("E8").Value = 0.3582277254281489
$ws.Range("F8").Value = 0.3748117494910135
$ws.Range("G8").Value = 9

# Row 9 (previously row 8's values)
$ws.Range("B9").Value = -0.02790046359007027
$ws.Range("C9").Value = 0.3120137538489823
$ws.Range("D9").Value = 0.1287774372279015
$ws.Range("E9").Value = 0.3588557331684997
$ws.Range("F9").Value = 0.3824716629792676
$ws.Range("G9").Value = 8

# Row 10 (previously row 9's values)
$ws.Range("B10").Value = -0.03188629724616485
$ws.Range("C10").Value = 0.2657389849834738
$ws.Range("D10").Value = 0.1191175486163901
$ws.Range("E10").Value = 0.3451341023665874
$ws.Range("F10").Value = 0.3711930514458828
$ws.Range("G10").Value = 7

# Row 11 (previously row 10's values)
$ws.Range("B11").Value = -0.0881379738596985
$ws.Range("C11").Value = 0.4006321887415465
$ws.Range("D11").Value = 0.2481252247171507
$ws.Range("E11").Value = 0.4981216966938408
$ws.Range("F11").Value = 0.5370552175868598
$ws.Range("G11").Value = 6
